$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.905.94"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.546.52"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "2.544.40"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "3.001.89"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "62.867.76"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "2.535.84"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "334.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.46"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.52%  "
$ws.Range("D31").Value = "0.0₃0811"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.65"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "404.56"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.14"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.400"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.75"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.75"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.600"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0241"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.92%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0965"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.13%  "

Write-Host "Updated 94 cells"
